$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: rename measurement column for the new event (NAUSEA -> FLUSHING)
$ws.Range("C1").Value = "Μετρήσεις για FLUSHING"

# Rows 2-51: refreshed drug list + updated counts/PRR values
$ws.Range("B2").Value = "NIACIN"
$ws.Range("C2").Value = 11032.0
$ws.Range("D2").Value = 53.45
$ws.Range("B3").Value = "DIMETHYL FUMARATE"
$ws.Range("C3").Value = 13832.0
$ws.Range("D3").Value = 23.6
$ws.Range("B4").Value = "DALFAMPRIDINE"
$ws.Range("C4").Value = 1780.0
$ws.Range("D4").Value = 6.27
$ws.Range("B5").Value = "FISH OIL"
$ws.Range("C5").Value = 2537.0
$ws.Range("D5").Value = 6.0
$ws.Range("B6").Value = "ASPIRIN 325 MG"
$ws.Range("C6").Value = 10709.0
$ws.Range("D6").Value = 4.95
$ws.Range("B7").Value = "ASPIRIN 81MG"
$ws.Range("C7").Value = 10707.0
$ws.Range("D7").Value = 4.95
$ws.Range("B8").Value = "ASPIRIN 81 MG"
$ws.Range("C8").Value = 10709.0
$ws.Range("D8").Value = 4.94
$ws.Range("B9").Value = "TREPROSTINIL"
$ws.Range("C9").Value = 1411.0
$ws.Range("D9").Value = 4.81
$ws.Range("B10").Value = "ASPIRIN"
$ws.Range("C10").Value = 10847.0
$ws.Range("D10").Value = 4.75
$ws.Range("B11").Value = "BENADRYL"
$ws.Range("C11").Value = 1436.0
$ws.Range("D11").Value = 4.31
$ws.Range("B12").Value = "TADALAFIL"
$ws.Range("C12").Value = 1826.0
$ws.Range("D12").Value = 3.79
$ws.Range("B13").Value = "GLATIRAMER ACETATE"
$ws.Range("C13").Value = 1094.0
$ws.Range("D13").Value = 3.24
$ws.Range("B14").Value = "DIPHENHYDRAMINE HYDROCHLORIDE"
$ws.Range("C14").Value = 1880.0
$ws.Range("D14").Value = 3.2
$ws.Range("B15").Value = "SILDENAFIL CITRATE"
$ws.Range("C15").Value = 1374.0
$ws.Range("D15").Value = 3.0
$ws.Range("B16").Value = "AMBRISENTAN"
$ws.Range("C16").Value = 1640.0
$ws.Range("D16").Value = 2.76
$ws.Range("B17").Value = "BACLOFEN"
$ws.Range("C17").Value = 1056.0
$ws.Range("D17").Value = 2.49
$ws.Range("B18").Value = "ROSUVASTATIN CALCIUM"
$ws.Range("C18").Value = 2338.0
$ws.Range("D18").Value = 2.47
$ws.Range("B19").Value = "ERGOCALCIFEROL"
$ws.Range("C19").Value = 3088.0
$ws.Range("D19").Value = 2.09
$ws.Range("B20").Value = "LISINOPRIL"
$ws.Range("C20").Value = 3416.0
$ws.Range("D20").Value = 2.04
$ws.Range("B21").Value = "WARFARIN SODIUM"
$ws.Range("C21").Value = 1269.0
$ws.Range("D21").Value = 1.86
$ws.Range("B22").Value = "ATORVASTATIN CALCIUM"
$ws.Range("C22").Value = 2037.0
$ws.Range("D22").Value = 1.85
$ws.Range("B23").Value = "CLOPIDOGREL"
$ws.Range("C23").Value = 2028.0
$ws.Range("D23").Value = 1.83
$ws.Range("B24").Value = "SIMVASTATIN"
$ws.Range("C24").Value = 2781.0
$ws.Range("D24").Value = 1.83
$ws.Range("B25").Value = "METOPROLOL"
$ws.Range("C25").Value = 2031.0
$ws.Range("D25").Value = 1.78
$ws.Range("B26").Value = "METOPROLOL TARTRATE"
$ws.Range("C26").Value = 2202.0
$ws.Range("D26").Value = 1.75
$ws.Range("B27").Value = "CLOPIDOGREL BISULFATE"
$ws.Range("C27").Value = 2048.0
$ws.Range("D27").Value = 1.74
$ws.Range("B28").Value = "HYDROCHLOROTHIAZIDE"
$ws.Range("C28").Value = 1263.0
$ws.Range("D28").Value = 1.74
$ws.Range("B29").Value = "LEVOTHYROXINE SODIUM"
$ws.Range("C29").Value = 2671.0
$ws.Range("D29").Value = 1.7
$ws.Range("B30").Value = "INFLIXIMAB"
$ws.Range("C30").Value = 1627.0
$ws.Range("D30").Value = 1.66
$ws.Range("B31").Value = "ATENOLOL"
$ws.Range("C31").Value = 1298.0
$ws.Range("D31").Value = 1.56
$ws.Range("B32").Value = "RANITIDINE HYDROCHLORIDE"
$ws.Range("C32").Value = 1179.0
$ws.Range("D32").Value = 1.56
$ws.Range("B33").Value = "IBUPROFEN"
$ws.Range("C33").Value = 1452.0
$ws.Range("D33").Value = 1.29
$ws.Range("B34").Value = "IBUPROFEN TABLETS"
$ws.Range("C34").Value = 1058.0
$ws.Range("D34").Value = 1.2
$ws.Range("B35").Value = "IBUPFROFEN"
$ws.Range("C35").Value = 1025.0
$ws.Range("D35").Value = 1.17
$ws.Range("B36").Value = "IBUPROFEN 200 MG"
$ws.Range("C36").Value = 1025.0
$ws.Range("D36").Value = 1.17
$ws.Range("B37").Value = "IBUPROFEN 200MG"
$ws.Range("C37").Value = 1025.0
$ws.Range("D37").Value = 1.17
$ws.Range("B38").Value = "IBUPROFEN ORAL"
$ws.Range("C38").Value = 1025.0
$ws.Range("D38").Value = 1.17
$ws.Range("B39").Value = "METFORMIN"
$ws.Range("C39").Value = 1420.0
$ws.Range("D39").Value = 1.13
$ws.Range("B40").Value = "ACETAMINOPHEN"
$ws.Range("C40").Value = 1559.0
$ws.Range("D40").Value = 1.08
$ws.Range("B41").Value = "GABAPENTIN"
$ws.Range("C41").Value = 1767.0
$ws.Range("D41").Value = 1.08
$ws.Range("B42").Value = "METFORMIN ER 500 MG"
$ws.Range("C42").Value = 1164.0
$ws.Range("D42").Value = 1.07
$ws.Range("B43").Value = "METFORMIN ER 750 MG"
$ws.Range("C43").Value = 1164.0
$ws.Range("D43").Value = 1.07
$ws.Range("B44").Value = "OMEPRAZOLE MAGNESIUM"
$ws.Range("C44").Value = 2482.0
$ws.Range("D44").Value = 1.07
$ws.Range("B45").Value = "AMLODIPINE"
$ws.Range("C45").Value = 1208.0
$ws.Range("D45").Value = 1.04
$ws.Range("B46").Value = "OMEPRAZOLE"
$ws.Range("C46").Value = 1567.0
$ws.Range("D46").Value = 0.94
$ws.Range("B47").Value = "ESOMEPRAZOLE MAGNESIUM"
$ws.Range("C47").Value = 1169.0
$ws.Range("D47").Value = 0.92
$ws.Range("B48").Value = "ESOMEPRAZOLE SODIUM"
$ws.Range("C48").Value = 1098.0
$ws.Range("D48").Value = 0.92
$ws.Range("B49").Value = "DEXAMETHASONE"
$ws.Range("C49").Value = 1095.0
$ws.Range("D49").Value = 0.84
$ws.Range("B50").Value = "FUROSEMIDE"
$ws.Range("C50").Value = 1665.0
$ws.Range("D50").Value = 0.83
$ws.Range("B51").Value = "PREDNISONE"
$ws.Range("C51").Value = 1240.0
$ws.Range("D51").Value = 0.63
